# Localization status report refresh ("Generate Report for Archive"):
# the handoff/translation status moved from "Ready for handoff" to
# "In Translation" everywhere it appears (Overview!E:F and the per-locale
# Status column on the zh-cn / de-de sheets), and the Status column was
# re-auto-fitted to the now-shorter text on those three sheets.

$wb = $excel.ActiveWorkbook

# 1) Update the status text on every worksheet. Cells.Replace only touches
#    exact/whole cell matches by default, so "Latest Handoff File" etc. are
#    left untouched.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# 2) Re-fit the columns that held the status text so they hug the shorter
#    string again (matches the narrower columns produced when the report is
#    regenerated with the new status value).
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = 12.5   # F: de-de status

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # C: Status

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # C: Status
